$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Q1_20_21" (first worksheet)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Insert two new blank columns at G:H, shifting the existing PVC/PVB columns
# (formerly G,H) to the right so they become I,J.
$ws1.Columns("G:H").Insert()

# New header row (row 2)
$ws1.Range("G2").Value = "VfM Category lower range"
$ws1.Range("H2").Value = "VfM Category upper range"
$ws1.Range("K2").Value = "Benefits Narrative"

# Row 3 - Mars
$ws1.Range("G3").Value = "Very High"
$ws1.Range("H3").Value = "Very High"
$ws1.Range("K3").Value = "All you need is love, love is all you need "

# Row 4 - SoT
$ws1.Range("G4").Value = "High"
$ws1.Range("H4").Value = "High"

# Row 6 - F9
$ws1.Range("G6").Value = "N/A"
$ws1.Range("H6").Value = "N/A"

# ---------------------------------------------------------------------------
# Sheet "Q4_19_20" (second worksheet)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Insert two new blank columns at G:H, shifting the existing PVC/PVB columns
# (formerly G,H) to the right so they become I,J.
$ws2.Columns("G:H").Insert()

# New header row (row 2)
$ws2.Range("G2").Value = "VfM Category lower range"
$ws2.Range("H2").Value = "VfM Category upper range"
$ws2.Range("K2").Value = "Benefits Narrative"

# Row 3 - Mars
$ws2.Range("K3").Value = "Hello is it me you’re looking for"

# Row 4 - SoT
$ws2.Range("K4").Value = "Please allow me to introduce myself I’m a man of wealth and taste."
